$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.355.73"
Set-TextValue $ws.Range("E2") "  -4.16%  "
Set-TextValue $ws.Range("D3") "1.862.06"
Set-TextValue $ws.Range("E3") "  -5.10%  "
Set-TextValue $ws.Range("D4") "0.9998"
Set-TextValue $ws.Range("E4") "  -1.29%  "
Set-TextValue $ws.Range("D5") "323.63"
Set-TextValue $ws.Range("E5") "  +0.25%  "
Set-TextValue $ws.Range("E6") "  -1.08%  "
Set-TextValue $ws.Range("E7") "  -5.90%  "
Set-TextValue $ws.Range("D8") "0.3868"
Set-TextValue $ws.Range("E8") "  -4.98%  "
Set-TextValue $ws.Range("D9") "47.91"
Set-TextValue $ws.Range("E9") "  -11.67%  "
Set-TextValue $ws.Range("D10") "0.07910"
Set-TextValue $ws.Range("E10") "  -6.78%  "
Set-TextValue $ws.Range("D11") "1.021"
Set-TextValue $ws.Range("E11") "  -3.78%  "
Set-TextValue $ws.Range("E12") "  -4.42%  "
Set-TextValue $ws.Range("D13") "1.859.10"
Set-TextValue $ws.Range("E13") "  -9.55%  "
Set-TextValue $ws.Range("D14") "5.901"
Set-TextValue $ws.Range("E14") "  -4.59%  "
Set-TextValue $ws.Range("D15") "7.154"
Set-TextValue $ws.Range("E15") "  -5.74%  "
Set-TextValue $ws.Range("D16") "0.9995"
Set-TextValue $ws.Range("E16") "  -1.37%  "
Set-TextValue $ws.Range("D17") "0.00001033"
Set-TextValue $ws.Range("E17") "  -3.77%  "
Set-TextValue $ws.Range("D18") "85.90"
Set-TextValue $ws.Range("E18") "  -5.66%  "
Set-TextValue $ws.Range("E19") "  -1.87%  "
Set-TextValue $ws.Range("D20") "17.14"
Set-TextValue $ws.Range("E20") "  -7.53%  "
Set-TextValue $ws.Range("D21") "1.001"
Set-TextValue $ws.Range("E21") "  -1.10%  "
Set-TextValue $ws.Range("D22") "5.525"
Set-TextValue $ws.Range("E22") "  -6.11%  "
Set-TextValue $ws.Range("D23") "27.359.31"
Set-TextValue $ws.Range("E23") "  -4.33%  "
Set-TextValue $ws.Range("D24") "10.84"
Set-TextValue $ws.Range("E24") "  -5.69%  "
Set-TextValue $ws.Range("D25") "2.267"
Set-TextValue $ws.Range("E25") "  -1.56%  "
Set-TextValue $ws.Range("D26") "2.088.61"
Set-TextValue $ws.Range("E26") "  -8.56%  "
Set-TextValue $ws.Range("E27") "  -2.48%  "
Set-TextValue $ws.Range("D28") "19.75"
Set-TextValue $ws.Range("E28") "  -2.90%  "
Set-TextValue $ws.Range("D29") "2.068"
Set-TextValue $ws.Range("E29") "  -5.36%  "
Set-TextValue $ws.Range("D30") "5.514"
Set-TextValue $ws.Range("E30") "  -6.25%  "
Set-TextValue $ws.Range("D31") "120.62"
Set-TextValue $ws.Range("E31") "  -3.45%  "
Set-TextValue $ws.Range("E32") "  +2.44%  "
Set-TextValue $ws.Range("D33") "0.09324"
Set-TextValue $ws.Range("E33") "  -3.74%  "
Set-TextValue $ws.Range("D34") "0.9377"
Set-TextValue $ws.Range("E34") "  -5.56%  "
Set-TextValue $ws.Range("E35") "  -2.56%  "
Set-TextValue $ws.Range("D36") "5.286"
Set-TextValue $ws.Range("E36") "  -6.43%  "
Set-TextValue $ws.Range("D37") "0.02238"
Set-TextValue $ws.Range("E37") "  -4.23%  "
Set-TextValue $ws.Range("D38") "0.06011"
Set-TextValue $ws.Range("E38") "  -3.76%  "
Set-TextValue $ws.Range("E39") "  -2.60%  "
Set-TextValue $ws.Range("D40") "8.272"
Set-TextValue $ws.Range("E40") "  -9.35%  "
Set-TextValue $ws.Range("D41") "0.9995"
Set-TextValue $ws.Range("E41") "  -1.11%  "
Set-TextValue $ws.Range("D42") "0.5920"
Set-TextValue $ws.Range("E42") "  -5.21%  "
Set-TextValue $ws.Range("D43") "0.1890"
Set-TextValue $ws.Range("E43") "  -1.54%  "
Set-TextValue $ws.Range("D44") "10.19"
Set-TextValue $ws.Range("E44") "  -9.44%  "
Set-TextValue $ws.Range("D45") "1.273"
Set-TextValue $ws.Range("E45") "  -6.30%  "
Set-TextValue $ws.Range("D46") "0.5646"
Set-TextValue $ws.Range("E46") "  -5.39%  "
Set-TextValue $ws.Range("D47") "12.00"
Set-TextValue $ws.Range("E47") "  -8.08%  "
Set-TextValue $ws.Range("B48") "NEARProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.928"
Set-TextValue $ws.Range("E48") "  -6.79%  "
Set-TextValue $ws.Range("B49") "PancakeSwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D49") "3.367"
Set-TextValue $ws.Range("E49") "  -1.42%  "
Set-TextValue $ws.Range("D50") "0.06796"
Set-TextValue $ws.Range("E50") "  -0.76%  "
Set-TextValue $ws.Range("D51") "108.00"
Set-TextValue $ws.Range("E51") "  -3.01%  "
